$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date from 45664 to 45665 for rows 2-37
for ($r = 2; $r -le 37; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45664) {
        $cell.Value = 45665
    }
}

# Swap row 36 and row 37 data for columns A (Beteckning) and G (Area (ha))
$ws.Cells.Item(36, 1).Value = "A 60500-2024"
$ws.Cells.Item(36, 7).Value = 0.8

$ws.Cells.Item(37, 1).Value = "A 60501-2024"
$ws.Cells.Item(37, 7).Value = 0.6
